$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new day's data as row 59 (A:D), matching the existing
# table's layout (date, weekday, hour, ranking).
# The leading apostrophe forces the date-looking string to be stored
# as literal text (same as the existing rows) instead of being
# auto-converted to a date serial number; re-applying the "Normal"
# style afterwards drops the quote-prefix formatting so the cell ends
# up styleless, just like its neighbours.
$ws.Range("A59").Value = "'2025/10/04"
$ws.Range("A59").Style = "Normal"
$ws.Range("B59").Value = "土"
$ws.Range("C59").Value = 16
$ws.Range("D59").Value = 43
